{"js": "const body = context.document.body;\n\n// Update the date heading paragraph (first paragraph in the body).\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items\");\nawait context.sync();\nparagraphs.items[0].insertText(\"2025-09-20 Saturday\", \"Replace\");\n\n// Update every arithmetic-problem cell in the practice table, in document order\n// (row-major: row 0 col 0..4, row 1 col 0..4, ...). Several old values repeat\n// (e.g. \"88-41=\", \"99-24=\") but map to different new values, so cells must be\n// addressed by position, not by searching for old text.\nconst tables = body.tables;\ntables.load(\"items\");\nawait context.sync();\nconst table = tables.items[0];\n\nconst newValues = [\n  [\"31-20=\", \"52+41=\", \"3+4=\", \"90-65=\", \"81-4=\"],\n  [\"11+85=\", \"30+0=\", \"73+18=\", \"91-88=\", \"66+26=\"],\n  [\"3+90=\", \"1+14=\", \"21-8=\", \"98-47=\", \"99-16=\"],\n  [\"75-16=\", \"33+62=\", \"59+34=\", \"45-31=\", \"27+55=\"],\n  [\"6+16=\", \"75-53=\", \"54-43=\", \"14+40=\", \"59+18=\"],\n  [\"68-11=\", \"86-25=\", \"83-8=\", \"66-32=\", \"37+20=\"],\n  [\"14+81=\", \"32+29=\", \"18+52=\", \"50+22=\", \"0+55=\"],\n  [\"40-26=\", \"94-28=\", \"46+22=\", \"87+0=\", \"60+0=\"],\n  [\"20-4=\", \"92-31=\", \"27+4=\", \"96+3=\", \"19-18=\"],\n  [\"67-35=\", \"1+71=\", \"88+11=\", \"99-77=\", \"70-17=\"],\n  [\"22+76=\", \"38+50=\", \"85-76=\", \"11+42=\", \"84-20=\"],\n  [\"33+48=\", \"58+36=\", \"4+23=\", \"93-27=\", \"22+20=\"],\n  [\"51+32=\", \"64-17=\", \"54+45=\", \"31-19=\", \"89+10=\"],\n  [\"57-31=\", \"69-22=\", \"41-6=\", \"59+36=\", \"59-54=\"],\n  [\"99-80=\", \"53-32=\", \"12+55=\", \"25+51=\", \"5+48=\"],\n  [\"37+54=\", \"50+1=\", \"53+34=\", \"21+7=\", \"14+82=\"],\n  [\"21+48=\", \"96+0=\", \"70-18=\", \"24+30=\", \"17+81=\"],\n  [\"97-6=\", \"8+5=\", \"0+44=\", \"96+2=\", \"88-68=\"],\n  [\"9+77=\", \"5+94=\", \"16+45=\", \"59-55=\", \"66-20=\"],\n  [\"90-4=\", \"67-60=\", \"51+24=\", \"53-37=\", \"30-3=\"],\n];\n\nfor (let r = 0; r < newValues.length; r++) {\n  for (let c = 0; c < newValues[r].length; c++) {\n    table.getCell(r, c).value = newValues[r][c];\n  }\n}\n\nawait context.sync();", "ps1": "$d = $word.ActiveDocument\n\n# Update the date heading (first paragraph in the body).\n$d.Paragraphs(1).Range.Text = \"2025-09-20 Saturday\"\n\n# Update every arithmetic-problem cell in the practice table, addressed by its\n# (row, column) position (1-based), in document order: row 1 col 1..5, row 2 col 1..5, ...\n# Several old values repeat (e.g. \"88-41=\", \"99-24=\") but map to different new\n# values, so cells must be addressed positionally rather than matched by old text.\n$tbl = $d.Tables.Item(1)\n\n$newValues = @(\n  @(\"31-20=\", \"52+41=\", \"3+4=\", \"90-65=\", \"81-4=\"),\n  @(\"11+85=\", \"30+0=\", \"73+18=\", \"91-88=\", \"66+26=\"),\n  @(\"3+90=\", \"1+14=\", \"21-8=\", \"98-47=\", \"99-16=\"),\n  @(\"75-16=\", \"33+62=\", \"59+34=\", \"45-31=\", \"27+55=\"),\n  @(\"6+16=\", \"75-53=\", \"54-43=\", \"14+40=\", \"59+18=\"),\n  @(\"68-11=\", \"86-25=\", \"83-8=\", \"66-32=\", \"37+20=\"),\n  @(\"14+81=\", \"32+29=\", \"18+52=\", \"50+22=\", \"0+55=\"),\n  @(\"40-26=\", \"94-28=\", \"46+22=\", \"87+0=\", \"60+0=\"),\n  @(\"20-4=\", \"92-31=\", \"27+4=\", \"96+3=\", \"19-18=\"),\n  @(\"67-35=\", \"1+71=\", \"88+11=\", \"99-77=\", \"70-17=\"),\n  @(\"22+76=\", \"38+50=\", \"85-76=\", \"11+42=\", \"84-20=\"),\n  @(\"33+48=\", \"58+36=\", \"4+23=\", \"93-27=\", \"22+20=\"),\n  @(\"51+32=\", \"64-17=\", \"54+45=\", \"31-19=\", \"89+10=\"),\n  @(\"57-31=\", \"69-22=\", \"41-6=\", \"59+36=\", \"59-54=\"),\n  @(\"99-80=\", \"53-32=\", \"12+55=\", \"25+51=\", \"5+48=\"),\n  @(\"37+54=\", \"50+1=\", \"53+34=\", \"21+7=\", \"14+82=\"),\n  @(\"21+48=\", \"96+0=\", \"70-18=\", \"24+30=\", \"17+81=\"),\n  @(\"97-6=\", \"8+5=\", \"0+44=\", \"96+2=\", \"88-68=\"),\n  @(\"9+77=\", \"5+94=\", \"16+45=\", \"59-55=\", \"66-20=\"),\n  @(\"90-4=\", \"67-60=\", \"51+24=\", \"53-37=\", \"30-3=\"),\n)\n\nfor ($r = 0; $r -lt $newValues.Count; $r++) {\n  for ($c = 0; $c -lt $newValues[$r].Count; $c++) {\n    $tbl.Cell($r + 1, $c + 1).Range.Text = $newValues[$r][$c]\n  }\n}\n"}
